$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# B2 used to hold "Home"; replace it with the date-range label.
$ws.Range("B2").Value = "12/25/2019 - 12/25/2019"

# C2, D2 and E2 used to hold "Flights", "Hotels" and "Car Rentals";
# clear their contents while keeping the existing cell formatting.
$ws.Range("C2:E2").ClearContents()

# Move the active selection from E2 to B2.
$ws.Range("B2").Select()
